$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Replace the plain dates in column A (rows 36-41) with descriptive text
#    that also records the new 10/19/23 SKAT/Burden run dates.
# ---------------------------------------------------------------------------
$dateText1 = "10/6/2023 & 10/19/23 for SKAT and Burden"
$dateText2 = "10/9/2023 & 10/19/23 for SKAT and Burden"

$ws.Range("A36").Value = $dateText1
$ws.Range("A37").Value = $dateText1
$ws.Range("A38").Value = $dateText1

$ws.Range("A39").Value = $dateText2
$ws.Range("A40").Value = $dateText2
$ws.Range("A41").Value = $dateText2

# ---------------------------------------------------------------------------
# 2. Fill in the newly-computed SKAT / Burden results (columns W:AH) for the
#    100v99 (rows 36-38) and 100v80 (rows 39-41) comparisons. These cells
#    used to hold a literal "-" placeholder with centered/bordered styling;
#    now they hold real numbers. Rows 36-38 revert to the plain default
#    style, rows 39-41 pick up a new right-aligned wrapped style.
# ---------------------------------------------------------------------------

# Row 36
$ws.Range("W36:AH36").Style = "Normal"
$ws.Range("W36").Value = 0
$ws.Range("X36").Value = 0
$ws.Range("Y36").Value = 0
$ws.Range("Z36").Value = 0.05
$ws.Range("AA36").Value = 0.04
$ws.Range("AB36").Value = 0.04
$ws.Range("AC36").Value = 0
$ws.Range("AD36").Value = 0.01
$ws.Range("AE36").Value = 0
$ws.Range("AF36").Value = 0.03
$ws.Range("AG36").Value = 0.07
$ws.Range("AH36").Value = 0.04

# Row 37
$ws.Range("W37:AH37").Style = "Normal"
$ws.Range("W37").Value = 0
$ws.Range("X37").Value = 0.85
$ws.Range("Y37").Value = 0.58
$ws.Range("Z37").Value = 0.05
$ws.Range("AA37").Value = 0.97
$ws.Range("AB37").Value = 0.82
$ws.Range("AC37").Value = 0
$ws.Range("AD37").Value = 0.7
$ws.Range("AE37").Value = 0.39
$ws.Range("AF37").Value = 0.03
$ws.Range("AG37").Value = 0.7
$ws.Range("AH37").Value = 0.52

# Row 38
$ws.Range("W38:AH38").Style = "Normal"
$ws.Range("W38").Value = 0
$ws.Range("X38").Value = 0
$ws.Range("Y38").Value = 0
$ws.Range("Z38").Value = 0.02
$ws.Range("AA38").Value = 0.04
$ws.Range("AB38").Value = 0.04
$ws.Range("AC38").Value = 0.01
$ws.Range("AD38").Value = 0
$ws.Range("AE38").Value = 0.02
$ws.Range("AF38").Value = 0.02
$ws.Range("AG38").Value = 0.03
$ws.Range("AH38").Value = 0.02

# Build the new right-aligned / wrapped style once on a scratch cell far off
# the used range, then paste just the formatting onto the target rows. Doing
# it this way (rather than setting HorizontalAlignment/VerticalAlignment
# directly on the destination range) produces a single clean style entry.
$scratch = $ws.Range("BZ1")
$scratch.HorizontalAlignment = -4152   # xlRight
$scratch.VerticalAlignment = -4107     # xlBottom (i.e. "unset")
$scratch.WrapText = $true

# Row 39
$scratch.Copy() | Out-Null
$ws.Range("W39:AH39").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("W39").Value = 0.01
$ws.Range("X39").Value = 0.02
$ws.Range("Y39").Value = 0.03
$ws.Range("Z39").Value = 0.04
$ws.Range("AA39").Value = 0.03
$ws.Range("AB39").Value = 0.07
$ws.Range("AC39").Value = 0
$ws.Range("AD39").Value = 0.02
$ws.Range("AE39").Value = 0.01
$ws.Range("AF39").Value = 0.03
$ws.Range("AG39").Value = 0
$ws.Range("AH39").Value = 0.02

# Row 40
$ws.Range("W40:AH40").Style = "Normal"
$ws.Range("W40").Value = 0.01
$ws.Range("X40").Value = 0.19
$ws.Range("Y40").Value = 0.07
$ws.Range("Z40").Value = 0.04
$ws.Range("AA40").Value = 0.12
$ws.Range("AB40").Value = 0.06
$ws.Range("AC40").Value = 0
$ws.Range("AD40").Value = 0.2
$ws.Range("AE40").Value = 0.1
$ws.Range("AF40").Value = 0.03
$ws.Range("AG40").Value = 0.08
$ws.Range("AH40").Value = 0.06

# Row 41
$ws.Range("W41:AH41").Style = "Normal"
$ws.Range("W41").Value = 0.01
$ws.Range("X41").Value = 0.02
$ws.Range("Y41").Value = 0.03
$ws.Range("Z41").Value = 0.03
$ws.Range("AA41").Value = 0.03
$ws.Range("AB41").Value = 0.07
$ws.Range("AC41").Value = 0
$ws.Range("AD41").Value = 0.02
$ws.Range("AE41").Value = 0.01
$ws.Range("AF41").Value = 0.02
$ws.Range("AG41").Value = 0.01
$ws.Range("AH41").Value = 0.03

$scratch.Clear() | Out-Null

# ---------------------------------------------------------------------------
# 3. Update the sheet's selection to reflect where the author left off after
#    entering the new data.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("G57").Select() | Out-Null
